$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.731.61"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "3.320.51"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "230.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "614.09"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("E9").Value = "  -0.03%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.946"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").Value = "3.318.64"
$ws.Range("E11").Value = "  -3.45%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "42.15"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.193"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.42%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.00"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").Value = "91.524.49"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").Value = "3.943.31"
$ws.Range("E16").Value = "  -3.35%  "
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("E18").Value = "  -3.49%  "
$ws.Range("D19").Value = "3.313.83"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -1.84%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.85"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.58%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.43"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +10.04%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "490.60"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.83%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.448"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -9.66%  "
$ws.Range("E25").Value = "  -2.79%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.10"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -6.88%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "89.68"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "11.78"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "3.498.03"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -5.62%  "
$ws.Range("E32").Value = "  +3.27%  "
$ws.Range("E33").Value = "  -4.41%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.170"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -5.47%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "28.07"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.81%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.524"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -6.35%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "553.09"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("E39").Value = "  -0.05%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "7.30"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("E42").Value = "  -4.97%  "
$ws.Range("E43").Value = "  -6.39%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E45").Value = "  +3.69%  "
$ws.Range("E46").Value = "  -0.40%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0409"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("E49").Value = "  -1.56%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "51.70"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -3.40%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.91"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
